$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "MiddleInitial" column (F) moves to the front of the table: cut it
# and insert it before column A, shifting Number/Gender/NameSet/Title/
# GivenName (A-E) one column to the right (B-F). Columns G onward are
# untouched.
$ws.Columns("F").Cut() | Out-Null
$ws.Columns("A").Insert() | Out-Null

# The sheet's defined name tracked the original data range starting at
# column A ($A$1:$AS$10). Now that column A holds the data that used to
# live in column B, the name should refer to the data starting at column B.
$name = $wb.Names.Item("FakeNameGenerator.com_e6cf4fc0")
$name.RefersTo = "=Tabelle1!`$B`$1:`$AS`$10"

# Reflect the editor's new selection: the whole of column G.
$ws.Columns("G").Select() | Out-Null

# Record the page setup that was added to the sheet (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
